$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 799.6667
$ws.Range("I18").Value = 799.6667
$ws.Range("K18").Value = 799.6667
$ws.Range("M18").Value = -515.6667

$ws.Range("H46").Value = 189373.86
$ws.Range("I46").Value = 267039
$ws.Range("J46").Value = 131125
$ws.Range("K46").Value = 801117
$ws.Range("L46").Value = 393375
$ws.Range("M46").Value = -800998
$ws.Range("N46").Value = -393613

$ws.Range("H60").Value = 189373.86
$ws.Range("I60").Value = 267039
$ws.Range("J60").Value = 131125
$ws.Range("K60").Value = 801117
$ws.Range("L60").Value = 393375
$ws.Range("M60").Value = -800633
$ws.Range("N60").Value = -394343

$ws.Range("H101").Value = 1506.1666
$ws.Range("I101").Value = 760.8
$ws.Range("J101").Value = 2038.5714
$ws.Range("K101").Value = 2282.4
$ws.Range("L101").Value = 6115.7142
$ws.Range("M101").Value = -660.3999999999996
$ws.Range("N101").Value = -9359.7142

$ws.Range("H137").Value = 2704446.2
$ws.Range("I137").Value = 4546898
$ws.Range("J137").Value = 2183.4
$ws.Range("K137").Value = 13640694
$ws.Range("L137").Value = 6550.200000000001
$ws.Range("M137").Value = -13638144
$ws.Range("N137").Value = -11650.2

$ws.Range("H138").Value = 1874636.8
$ws.Range("I138").Value = 1266.85
$ws.Range("J138").Value = 3403918.2
$ws.Range("K138").Value = 3800.55
$ws.Range("L138").Value = 10211754.6
$ws.Range("M138").Value = 1339.45
$ws.Range("N138").Value = -10222034.6

$ws.Range("H141").Value = 2991.9546
$ws.Range("I141").Value = 2420.1428
$ws.Range("J141").Value = 15000
$ws.Range("K141").Value = 7260.428400000001
$ws.Range("L141").Value = 45000
$ws.Range("M141").Value = -2080.428400000001
$ws.Range("N141").Value = -55360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1310.97
$ws.Range("I32").Value = 876.4286
$ws.Range("J32").Value = 3592.3125
$ws.Range("K32").Value = 876.4286
$ws.Range("L32").Value = 3592.3125
$ws.Range("M32").Value = -589.4286
$ws.Range("N32").Value = -4166.3125

$ws.Range("H135").Value = 48675
$ws.Range("J135").Value = 48675
$ws.Range("L135").Value = 48675
$ws.Range("N135").Value = -58815

$ws.Range("H139").Value = 44166
$ws.Range("J139").Value = 44166
$ws.Range("L139").Value = 44166
$ws.Range("N139").Value = -54446

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1445.9642
$ws.Range("I107").Value = 1220.6086
$ws.Range("J107").Value = 2482.6
$ws.Range("K107").Value = 1220.6086
$ws.Range("L107").Value = 2482.6
$ws.Range("M107").Value = 699.3914
$ws.Range("N107").Value = -6322.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 7384685
$ws.Range("I6").Value = 9600061
$ws.Range("J6").Value = 100
$ws.Range("K6").Value = 9600061
$ws.Range("L6").Value = 100
$ws.Range("M6").Value = -9599948
$ws.Range("N6").Value = -326

$ws.Range("H7").Value = 250.25
$ws.Range("I7").Value = 67.75
$ws.Range("J7").Value = 432.75
$ws.Range("K7").Value = 67.75
$ws.Range("L7").Value = 432.75
$ws.Range("M7").Value = 45.25
$ws.Range("N7").Value = -658.75

$ws.Range("H31").Value = 2403.6345
$ws.Range("I31").Value = 1154.8857
$ws.Range("J31").Value = 4974.5884
$ws.Range("K31").Value = 1154.8857
$ws.Range("L31").Value = 4974.5884
$ws.Range("M31").Value = -859.8857
$ws.Range("N31").Value = -5564.5884

$ws.Range("H34").Value = 2403.6345
$ws.Range("I34").Value = 1154.8857
$ws.Range("J34").Value = 4974.5884
$ws.Range("K34").Value = 1154.8857
$ws.Range("L34").Value = 4974.5884
$ws.Range("M34").Value = -952.8857
$ws.Range("N34").Value = -5378.5884

$ws.Range("H41").Value = 12000
$ws.Range("I41").Value = 6000
$ws.Range("K41").Value = 6000
$ws.Range("M41").Value = -5572

$ws.Range("H50").Value = 24673.6
$ws.Range("J50").Value = 24673.6
$ws.Range("L50").Value = 24673.6
$ws.Range("N50").Value = -25923.6

$ws.Range("H51").Value = 24824.25
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 24824.25
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 24824.25
$ws.Range("N51").Value = -26296.25
$ws.Range("M51").ClearContents()

$ws.Range("H59").Value = 32914.3
$ws.Range("J59").Value = 32914.3
$ws.Range("L59").Value = 32914.3
$ws.Range("N59").Value = -35204.3

$ws.Range("H60").Value = 9500
$ws.Range("I60").Value = 4000
$ws.Range("K60").Value = 4000
$ws.Range("M60").Value = -3489

$ws.Range("H61").Value = 24824.25
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 24824.25
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 24824.25
$ws.Range("N61").Value = -25520.25
$ws.Range("M61").ClearContents()

$ws.Range("H68").Value = 33636
$ws.Range("J68").Value = 38295
$ws.Range("L68").Value = 38295
$ws.Range("N68").Value = -39793

$ws.Range("H71").Value = 33636
$ws.Range("J71").Value = 38295
$ws.Range("L71").Value = 114885
$ws.Range("N71").Value = -122373

$ws.Range("H74").Value = 25300.348
$ws.Range("J74").Value = 25300.348
$ws.Range("L74").Value = 25300.348
$ws.Range("N74").Value = -27048.348

$ws.Range("H77").Value = 25300.348
$ws.Range("J77").Value = 25300.348
$ws.Range("L77").Value = 75901.04400000001
$ws.Range("N77").Value = -84637.04400000001

$ws.Range("H107").Value = 304.17392
$ws.Range("I107").Value = 324.09525
$ws.Range("J107").Value = 95
$ws.Range("K107").Value = 324.09525
$ws.Range("L107").Value = 95
$ws.Range("M107").Value = 1595.90475
$ws.Range("N107").Value = -3935

$ws.Range("H132").Value = 16327.53
$ws.Range("I132").Value = 1123.6792
$ws.Range("J132").Value = 78312.46000000001
$ws.Range("K132").Value = 3371.0376
$ws.Range("L132").Value = 234937.38
$ws.Range("M132").Value = -841.0376000000001
$ws.Range("N132").Value = -239997.38

$ws.Range("H134").Value = 18233.35
$ws.Range("I134").Value = 1038.34
$ws.Range("K134").Value = 3115.02
$ws.Range("M134").Value = -580.0199999999995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5884582.5
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 6252362.5
$ws.Range("K4").Value = 300
$ws.Range("L4").Value = 18757087.5
$ws.Range("M4").Value = -188
$ws.Range("N4").Value = -18757311.5

$ws.Range("H56").Value = 154940.14
$ws.Range("I56").Value = 154940.14
$ws.Range("K56").Value = 154940.14
$ws.Range("M56").Value = -154410.14

$ws.Range("H118").Value = 3210
$ws.Range("I118").Value = 1300
$ws.Range("J118").Value = 3687.5
$ws.Range("K118").Value = 3900
$ws.Range("L118").Value = 11062.5
$ws.Range("M118").Value = -2657
$ws.Range("N118").Value = -13548.5

$ws.Range("H119").Value = 3332.9
$ws.Range("I119").Value = 1221.5
$ws.Range("K119").Value = 3664.5
$ws.Range("M119").Value = 1173.5

$ws.Range("H138").Value = 3229.4736
$ws.Range("I138").Value = 2466
$ws.Range("J138").Value = 4077.7778
$ws.Range("K138").Value = 7398
$ws.Range("L138").Value = 12233.3334
$ws.Range("M138").Value = -2258
$ws.Range("N138").Value = -22513.3334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2250
$ws.Range("J126").Value = 3000
$ws.Range("L126").Value = 9000
$ws.Range("N126").Value = -13940

$ws.Range("H128").Value = 55000
$ws.Range("J128").Value = 55000
$ws.Range("L128").Value = 55000
$ws.Range("N128").Value = -64960

$ws.Range("H138").Value = 55000
$ws.Range("J138").Value = 55000
$ws.Range("L138").Value = 55000
$ws.Range("N138").Value = -65280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 75878.57000000001
$ws.Range("I136").Value = 59706.53
$ws.Range("J136").Value = 100871.73
$ws.Range("K136").Value = 179119.59
$ws.Range("L136").Value = 302615.19
$ws.Range("M136").Value = -176569.59
$ws.Range("N136").Value = -307715.19

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()

$ws.Range("H62").Value = 5000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 5000
$ws.Range("N62").Value = -6248
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 5000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 25000
$ws.Range("N65").Value = -31240
$ws.Range("M65").ClearContents()

$ws.Range("H132").Value = 39285.77
$ws.Range("I132").Value = 22514.674
$ws.Range("J132").Value = 167864.17
$ws.Range("K132").Value = 67544.022
$ws.Range("L132").Value = 503592.51
$ws.Range("M132").Value = -65014.022
$ws.Range("N132").Value = -508652.51
